$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows (2-7) and add new rows (8-10) to reflect updated TPM data
# Row 2
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Bmp10"
$ws.Range("C2").Value = "Bmpr1b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.178715
$ws.Range("H2").Value = 0.536145
$ws.Range("I2").Value = 0.09904930989061336
$ws.Range("J2").Value = 0.09904930989061336
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09074700000000001
$ws.Range("N2").Value = 0.272241
$ws.Range("O2").Value = 0.04501122713837679
$ws.Range("P2").Value = 0.04501122713837678
$ws.Range("Q2").Value = 0.016217850105
$ws.Range("R2").Value = 0.145960650945
$ws.Range("S2").Value = 0.004458330985385869
$ws.Range("T2").Value = 0.004458330985385868

# Row 3
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Bmp10"
$ws.Range("C3").Value = "Bmpr1b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.178715
$ws.Range("H3").Value = 0.536145
$ws.Range("I3").Value = 0.09904930989061336
$ws.Range("J3").Value = 0.09904930989061336
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.394935666666666
$ws.Range("N3").Value = 4.184806999999999
$ws.Range("O3").Value = 0.6918990835593063
$ws.Range("P3").Value = 0.6918990835593063
$ws.Range("Q3").Value = 0.2492959276683333
$ws.Range("R3").Value = 2.243663349015
$ws.Range("S3").Value = 0.06853212674049712
$ws.Range("T3").Value = 0.06853212674049712

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Bmp10"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.178715
$ws.Range("H4").Value = 0.536145
$ws.Range("I4").Value = 0.09904930989061336
$ws.Range("J4").Value = 0.09904930989061336
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5304143333333333
$ws.Range("N4").Value = 1.591243
$ws.Range("O4").Value = 0.263089689302317
$ws.Range("P4").Value = 0.263089689302317
$ws.Range("Q4").Value = 0.09479299758166666
$ws.Range("R4").Value = 0.8531369782349999
$ws.Range("S4").Value = 0.02605885216473039
$ws.Range("T4").Value = 0.02605885216473039

# Row 5
$ws.Range("A5").Value = "Neutrophils"
$ws.Range("B5").Value = "Bmp10"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.450498333333333
$ws.Range("H5").Value = 4.351495
$ws.Range("I5").Value = 0.8039104659046613
$ws.Range("J5").Value = 0.8039104659046612
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.09074700000000001
$ws.Range("N5").Value = 0.272241
$ws.Range("O5").Value = 0.04501122713837679
$ws.Range("P5").Value = 0.04501122713837678
$ws.Range("Q5").Value = 0.131628372255
$ws.Range("R5").Value = 1.184655350295
$ws.Range("S5").Value = 0.03618499657975302
$ws.Range("T5").Value = 0.03618499657975301

# Row 6
$ws.Range("A6").Value = "Neutrophils"
$ws.Range("B6").Value = "Bmp10"
$ws.Range("C6").Value = "Bmpr1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.450498333333333
$ws.Range("H6").Value = 4.351495
$ws.Range("I6").Value = 0.8039104659046613
$ws.Range("J6").Value = 0.8039104659046612
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.394935666666666
$ws.Range("N6").Value = 4.184806999999999
$ws.Range("O6").Value = 0.6918990835593063
$ws.Range("P6").Value = 0.6918990835593063
$ws.Range("Q6").Value = 2.023351859607222
$ws.Range("R6").Value = 18.21016673646499
$ws.Range("S6").Value = 0.5562249146231701
$ws.Range("T6").Value = 0.55622491462317

# Row 7
$ws.Range("A7").Value = "Neutrophils"
$ws.Range("B7").Value = "Bmp10"
$ws.Range("C7").Value = "Bmpr1b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.450498333333333
$ws.Range("H7").Value = 4.351495
$ws.Range("I7").Value = 0.8039104659046613
$ws.Range("J7").Value = 0.8039104659046612
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5304143333333333
$ws.Range("N7").Value = 1.591243
$ws.Range("O7").Value = 0.263089689302317
$ws.Range("P7").Value = 0.263089689302317
$ws.Range("Q7").Value = 0.769365106476111
$ws.Range("R7").Value = 6.924285958285
$ws.Range("S7").Value = 0.2115005547017383
$ws.Range("T7").Value = 0.2115005547017382

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Bmp10"
$ws.Range("C8").Value = "Bmpr1b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.17509
$ws.Range("H8").Value = 0.52527
$ws.Range("I8").Value = 0.09704022420472538
$ws.Range("J8").Value = 0.09704022420472537
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.09074700000000001
$ws.Range("N8").Value = 0.272241
$ws.Range("O8").Value = 0.04501122713837679
$ws.Range("P8").Value = 0.04501122713837678
$ws.Range("Q8").Value = 0.01588889223
$ws.Range("R8").Value = 0.14300003007
$ws.Range("S8").Value = 0.004367899573237903
$ws.Range("T8").Value = 0.004367899573237902

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Bmp10"
$ws.Range("C9").Value = "Bmpr1b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.17509
$ws.Range("H9").Value = 0.52527
$ws.Range("I9").Value = 0.09704022420472538
$ws.Range("J9").Value = 0.09704022420472537
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.394935666666666
$ws.Range("N9").Value = 4.184806999999999
$ws.Range("O9").Value = 0.6918990835593063
$ws.Range("P9").Value = 0.6918990835593063
$ws.Range("Q9").Value = 0.2442392858766666
$ws.Range("R9").Value = 2.198153572889999
$ws.Range("S9").Value = 0.0671420421956391
$ws.Range("T9").Value = 0.0671420421956391

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Bmp10"
$ws.Range("C10").Value = "Bmpr1b"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.17509
$ws.Range("H10").Value = 0.52527
$ws.Range("I10").Value = 0.09704022420472538
$ws.Range("J10").Value = 0.09704022420472537
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.5304143333333333
$ws.Range("N10").Value = 1.591243
$ws.Range("O10").Value = 0.263089689302317
$ws.Range("P10").Value = 0.263089689302317
$ws.Range("Q10").Value = 0.09287024562333333
$ws.Range("R10").Value = 0.83583221061
$ws.Range("S10").Value = 0.02553028243584839
$ws.Range("T10").Value = 0.02553028243584838
